$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new team-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the existing header row (A1:AC1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record (Wins/Losses/Ties) for every player row
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 91
    $ws.Cells.Item($r, 31).Value = 71
    $ws.Cells.Item($r, 32).Value = 0
}
